$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1 - copy style from H1 (bold/centered/bordered header style)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data values for column I (I0)
$iValues = @(7, 8, 9, 5, 5, 5, 7, 6, 5, 9, 3, 1)
# Data values for column J (IF)
$jValues = @(7, 8, 9, 5, 5, 6, 7, 6, 5, 9, 3, 1)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
